$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings (e.g. "0.999")
# are stored as literal text, matching the source data (inline strings), not
# auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "51.783.15"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "2.931.25"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "352.62"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("D6").Value = "107.05"
$ws.Range("E6").Value = "  -6.01%  "
$ws.Range("D7").Value = "0.556"
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.614"
$ws.Range("E9").Value = "  -1.47%  "
$ws.Range("D10").Value = "37.84"
$ws.Range("E10").Value = "  -4.59%  "
$ws.Range("D11").Value = "0.138"
$ws.Range("E11").Value = "  +0.97%  "
$ws.Range("D12").Value = "0.0854"
$ws.Range("E12").Value = "  -2.60%  "
$ws.Range("D13").Value = "19.04"
$ws.Range("E13").Value = "  -3.90%  "
$ws.Range("D14").Value = "3.387.63"
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("D15").Value = "7.58"
$ws.Range("E15").Value = "  -2.03%  "
$ws.Range("D16").Value = "2.923.89"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").Value = "0.967"
$ws.Range("E17").Value = "  -1.97%  "
$ws.Range("D18").Value = "51.677.56"
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("D19").Value = "3.46"
$ws.Range("E19").Value = "  +3.73%  "
$ws.Range("D20").Value = "7.36"
$ws.Range("E20").Value = "  -3.18%  "
$ws.Range("D21").Value = "13.46"
$ws.Range("E21").Value = "  -4.50%  "
$ws.Range("D22").Value = "0.0₃0965"
$ws.Range("E22").Value = "  -1.81%  "
$ws.Range("D23").Value = "69.04"
$ws.Range("E23").Value = "  -3.02%  "
$ws.Range("D24").Value = "262.19"
$ws.Range("E24").Value = "  -2.75%  "
$ws.Range("D25").Value = "2.73"
$ws.Range("E25").Value = "  -3.03%  "
$ws.Range("E26").Value = "  -6.18%  "
$ws.Range("D27").Value = "26.56"
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").Value = "7.50"
$ws.Range("E28").Value = "  +9.39%  "
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("D30").Value = "0.104"
$ws.Range("E30").Value = "  +0.86%  "
$ws.Range("D31").Value = "10.26"
$ws.Range("E31").Value = "  -3.94%  "
$ws.Range("D32").Value = "35.65"
$ws.Range("E32").Value = "  -4.98%  "
$ws.Range("E33").Value = "  -4.98%  "
$ws.Range("D34").Value = "5.95"
$ws.Range("E34").Value = "  -1.56%  "
$ws.Range("D35").Value = "51.15"
$ws.Range("E35").Value = "  -3.79%  "
$ws.Range("D36").Value = "0.0429"
$ws.Range("E36").Value = "  -5.36%  "
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("D38").Value = "3.15"
$ws.Range("E38").Value = "  -5.83%  "
$ws.Range("D39").Value = "1.97"
$ws.Range("E39").Value = "  -3.81%  "
$ws.Range("D40").Value = "17.73"
$ws.Range("E40").Value = "  -5.74%  "
$ws.Range("D41").Value = "2.67"
$ws.Range("E41").Value = "  -1.66%  "
$ws.Range("E42").Value = "  -0.91%  "
$ws.Range("D43").Value = "22.51"
$ws.Range("E43").Value = "  -2.83%  "
$ws.Range("D44").Value = "119.70"
$ws.Range("E44").Value = "  +1.54%  "
$ws.Range("E45").Value = "  -1.40%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.100.88"
$ws.Range("E46").Value = "  -3.63%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "3.33"
$ws.Range("E47").Value = "  -5.86%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "2.36"
$ws.Range("E48").Value = "  -6.11%  "
$ws.Range("D49").Value = "3.214.15"
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("D50").Value = "0.239"
$ws.Range("E50").Value = "  -6.00%  "
$ws.Range("D51").Value = "0.0341"
$ws.Range("E51").Value = "  -4.04%  "

# Restore default (Normal) style on column D so no stray per-cell style
# references remain (keeps cellXfs/style indices matching the source).
$ws.Range("D2:D51").Style = "Normal"

